$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write brand-new text values in the exact order they were first
#     introduced so the shared-string table is appended in the same
#     sequence as the target file. ---

# 1) "Akif Elcan" - new row 159
$ws.Range("A159").Value = "Akif Elcan"

# 2) "diyarbakır telekom" - first used on row 5
$ws.Range("B5").Value  = "diyarbakır telekom"
$ws.Range("B6").Value  = "diyarbakır telekom"
$ws.Range("B12").Value = "diyarbakır telekom"
$ws.Range("B159").Value = "diyarbakır telekom"

# 3) "Anayurt Komşu" - row 10
$ws.Range("B10").Value = "Anayurt Komşu"

# 4) "Şerefiye" - row 26 onward
$ws.Range("B26").Value = "Şerefiye"
$ws.Range("B56").Value = "Şerefiye"
$ws.Range("B71").Value = "Şerefiye"
$ws.Range("B72").Value = "Şerefiye"
$ws.Range("B76").Value = "Şerefiye"

# 5) "Yusuf-Kaniye Özsoy" - corrected name, row 53
$ws.Range("A53").Value = "Yusuf-Kaniye Özsoy"

# 6) "Özgür Demir" - new row 160
$ws.Range("A160").Value = "Özgür Demir"

# 7) "Caner Kar" - new row 161
$ws.Range("A161").Value = "Caner Kar"

# --- Remaining values reuse existing shared strings, order irrelevant ---
$ws.Range("B2").Value  = "Komşu"
$ws.Range("B7").Value  = "Komşu"
$ws.Range("C159").Value = "Çeyrek"
$ws.Range("C160").Value = "1 gr altın"
$ws.Range("C161").Value = "Çeyrek"

# --- Restore the selection to I17 (matches the saved view in the target file) ---
$ws.Range("I17").Select()
